$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows below the data table ---

# Row 12: average of column J (|S*|/n) across the 10 data rows, bold
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# Rows 14-17: labeled summary statistics
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Style the B14:B17 summary values: bold, size 12, vertically centered, taller rows
$rngB = $ws.Range("B14:B17")
$rngB.Font.Bold = $true
$rngB.Font.Size = 12
$rngB.VerticalAlignment = -4108
$ws.Range("A14:B17").RowHeight = 15.6

# --- Page setup (portrait, paper size 9 = A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Leave the active selection on J12, matching the final authoring state ---
[void]$ws.Range("J12").Select()
